# Fruta / hortaliza, semanal
# The data rows (2-12) got re-ordered (re-sorted by date after the week's
# worth of new data came in), which shows up in the OOXML diff as a set of
# per-cell value changes. Capture every row's full contents first, then
# write them back out in the new row order so every cell ends up correct.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row number -> original row number that its data came from.
$rowMap = @{
    2  = 11
    3  = 4
    4  = 5
    5  = 2
    6  = 8
    7  = 9
    8  = 12
    9  = 7
    10 = 10
    11 = 6
    12 = 3
}

$firstCol = 1   # A
$lastCol  = 18  # R

# Snapshot all the original values for rows 2-12, columns A-R.
$original = @{}
for ($r = 2; $r -le 12; $r++) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $original[$r] = $rowVals
}

# Write back values according to the new row mapping.
foreach ($newRow in $rowMap.Keys) {
    $srcRow = $rowMap[$newRow]
    $srcVals = $original[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($newRow, $c).Value2 = $srcVals[$c]
    }
}
